$d = $word.ActiveDocument

# The placeholder paragraph ("  ") is the final paragraph of the document body.
$target = $d.Paragraphs.Last
$r = $target.Range

# Confirm we located the intended placeholder text before replacing it.
$found = $r.Find.Execute("  ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $found) {
    throw "Could not locate the placeholder paragraph to replace."
}

$newContent = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Integrating the systems to communicate with the rest of the network went smoother than expected. Other than the new functions needed for it, there was only two changes that needed to be made. First, the alert function needed to call the new function that sends messages to the network. Second, the function calls in the main loop had to be moved to a scheduled task, as the </w:t></w:r><w:r><w:t>mesh system requires the main loop to only call the mesh’s update function. As much of the code that was previously in the main loop only served as a method to delay the device from producing alerts while the rolling average of its position calibrated, it could be blended into the setup function to remove an if statement that became unnecessary after the first second of the device running.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Running the function to </w:t></w:r><w:r><w:t>read</w:t></w:r><w:r><w:t xml:space="preserve"> the IMU on a scheduled task introduced a delay in how often it was </w:t></w:r><w:r><w:t>called</w:t></w:r><w:r><w:t xml:space="preserve">, which had the unexpected side effect of increasing the device’s ability to tell the difference between negligible movement from bumping the table and significant movement that would justify </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">generating an alert, as in the case of bumping the table, the device would return to its original position before </w:t></w:r><w:r><w:t>the next call of the reading function.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The device </w:t></w:r><w:r><w:t>could be improved by mounting it onto a permanent circuit board, rather than just a temporary breadboard. It also currently relies on being plugged into a usb port to receive power, which ideally should be replaced with a rechargeable battery.</w:t></w:r></w:p>
'

[void]$target.Range.InsertXML($newContent)
